$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Duplicate row 2 (A2:N2) into row 3: first copy values+formats, then re-copy
# just the formats on top so the per-cell style indexes survive the paste.
$ws.Range("A2:N2").Copy()
$ws.Range("A3:N3").PasteSpecial(-4104) # xlPasteAll
$excel.CutCopyMode = $false
$ws.Range("A2:N2").Copy()
$ws.Range("A3:N3").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Overwrite the cells that differ from row 2. A leading apostrophe forces a
# text literal (rather than a number) while keeping the existing quote-prefix
# cell style instead of minting a new one.
$ws.Range("A3").Value = "'2"
$ws.Range("B3").Value = "'423698529"
$ws.Range("D3").Value = "'423698529"
$ws.Range("E3").Value = "'1234"
$ws.Range("M3").Value = "'pruebasqa90"

# Add hyperlink on N3 pointing to the same mailto address as N2
$ws.Hyperlinks.Add($ws.Range("N3"), "mailto:jalzate@todo1.net")

# Adding the hyperlink can mint a fresh style for N3; reapply N2's format so
# N3 keeps sharing the original hyperlink cell style.
$ws.Range("N2").Copy()
$ws.Range("N3").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Update the active selection to B3
$ws.Range("B3").Select()
